$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "1.008") are preserved verbatim instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.506.53'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '1.878.24'
$ws.Range("E3").Value = '  -2.03%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  -3.29%  '
$ws.Range("D5").Value = '314.73'
$ws.Range("E5").Value = '  -2.26%  '
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  -2.71%  '
$ws.Range("D7").Value = '0.5104'
$ws.Range("E7").Value = '  -2.21%  '
$ws.Range("D8").Value = '0.3937'
$ws.Range("E8").Value = '  -0.92%  '
$ws.Range("D9").Value = '0.08407'
$ws.Range("E9").Value = '  +0.09%  '
$ws.Range("D10").Value = '1.110'
$ws.Range("E10").Value = '  -2.69%  '
$ws.Range("D11").Value = '41.60'
$ws.Range("E11").Value = '  -2.97%  '
$ws.Range("D12").Value = '6.251'
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("D13").Value = '1.871.56'
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").Value = '20.49'
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("D15").Value = '7.260'
$ws.Range("E15").Value = '  -1.20%  '
$ws.Range("E16").Value = '  -3.41%  '
$ws.Range("D17").Value = '0.00001105'
$ws.Range("E17").Value = '  -1.26%  '
$ws.Range("D18").Value = '90.93'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").Value = '0.06718'
$ws.Range("E19").Value = '  -2.10%  '
$ws.Range("D20").Value = '17.69'
$ws.Range("E20").Value = '  -2.13%  '
$ws.Range("D21").Value = '1.006'
$ws.Range("E21").Value = '  -2.86%  '
$ws.Range("D22").Value = '5.949'
$ws.Range("E22").Value = '  -2.96%  '
$ws.Range("D23").Value = '28.513.22'
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("D25").Value = '2.257'
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").Value = '2.095.44'
$ws.Range("E26").Value = '  -1.32%  '
$ws.Range("D27").Value = '161.17'
$ws.Range("E27").Value = '  -1.61%  '
$ws.Range("D28").Value = '20.70'
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("D29").Value = '2.380'
$ws.Range("E29").Value = '  -3.46%  '
$ws.Range("D30").Value = '126.34'
$ws.Range("E30").Value = '  -1.58%  '
$ws.Range("D31").Value = '0.1048'
$ws.Range("E31").Value = '  -2.40%  '
$ws.Range("D32").Value = '1.049'
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("D33").Value = '5.780'
$ws.Range("E33").Value = '  -3.69%  '
$ws.Range("D34").Value = '3.616'
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("D35").Value = '0.02446'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("D36").Value = '0.06509'
$ws.Range("E36").Value = '  -2.69%  '
$ws.Range("D37").Value = '0.2186'
$ws.Range("E37").Value = '  -2.01%  '
$ws.Range("D38").Value = '8.919'
$ws.Range("E38").Value = '  -5.84%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("D40").Value = '1.195'
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").Value = '5.080'
$ws.Range("E41").Value = '  +1.05%  '
$ws.Range("D42").Value = '0.6438'
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("D43").Value = '11.16'
$ws.Range("E43").Value = '  -0.82%  '
$ws.Range("D44").Value = '1.006'
$ws.Range("E44").Value = '  -2.80%  '
$ws.Range("D45").Value = '0.6068'
$ws.Range("E45").Value = '  -1.74%  '
$ws.Range("D46").Value = '13.01'
$ws.Range("E46").Value = '  -1.98%  '
$ws.Range("D47").Value = '3.696'
$ws.Range("D48").Value = '2.009'
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("D49").Value = '122.46'
$ws.Range("E49").Value = '  -0.95%  '
$ws.Range("D50").Value = '1.208'
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("D51").Value = '1.194'
$ws.Range("E51").Value = '  -9.17%  '

# Restore the default style on column D so the only lasting change is the
# cell text, matching the original (unstyled) cells.
$ws.Range("D2:D51").Style = "Normal"
